# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns
# to the latest scraped values, preserving the original text cell type
# (NumberFormat "@" forces a numeric-looking string to be stored as text,
# then Style is reset to "Normal" so no stray formatting is introduced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.884.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.180.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.89%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.86%  "

$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.79"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.33%  "

$ws.Range("E11").Value = "  -4.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.86%  "

$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.506.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.177.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.99%  "

$ws.Range("E18").Value = "  -6.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.731.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.91%  "

$ws.Range("E20").Value = "  -1.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.62%  "

$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.96%  "

$ws.Range("E28").Value = "  -9.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.98%  "

$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("E31").Value = "  -5.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0777"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("E35").Value = "  -6.51%  "

$ws.Range("E36").Value = "  -3.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.10%  "

$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "12.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.64%  "

$ws.Range("E41").Value = "  -2.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "59.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.191"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0968"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.27%  "

$ws.Range("E48").Value = "  -4.63%  "

$ws.Range("E49").Value = "  -4.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.32%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.417"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.55%  "
